# Update the cryptos worksheet with latest prices / volume percentages
# (and a couple of rank-order swaps / new listings) per the upstream
# GitHub Actions scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:E on data rows 2-51 are stored as literal text (not numbers),
# even though many of the Price values look numeric (e.g. "0.9992",
# "29.097.17"). Force text formatting before writing so Excel doesn't
# silently reinterpret them as numbers/dates.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# Auto-generated cell updates from the diff
$ws.Range('D2').Value = '29.097.17'
$ws.Range('E2').Value = '  -1.74%  '
$ws.Range('D3').Value = '1.836.78'
$ws.Range('E3').Value = '  -1.24%  '
$ws.Range('D4').Value = '0.9992'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '240.50'
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').Value = '0.6773'
$ws.Range('E6').Value = '  -2.61%  '
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('D9').Value = '0.07455'
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('D10').Value = '23.07'
$ws.Range('E10').Value = '  -2.84%  '
$ws.Range('D11').Value = '0.07663'
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').Value = '1.834.20'
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').Value = '5.022'
$ws.Range('E13').Value = '  -2.62%  '
$ws.Range('D14').Value = '0.6760'
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('D15').Value = '86.27'
$ws.Range('E15').Value = '  -6.24%  '
$ws.Range('D16').Value = '6.161'
$ws.Range('E16').Value = '  -6.28%  '
$ws.Range('D17').Value = '29.092.46'
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('D18').Value = '0.000008257'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').Value = '2.071.34'
$ws.Range('E19').Value = '  -1.53%  '
$ws.Range('D20').Value = '227.60'
$ws.Range('E20').Value = '  -5.49%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = '7.344'
$ws.Range('E23').Value = '  -3.41%  '
$ws.Range('D24').Value = '0.9999'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = '160.62'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('D26').Value = '0.1431'
$ws.Range('E26').Value = '  -4.63%  '
$ws.Range('D27').Value = '8.703'
$ws.Range('E27').Value = '  -2.52%  '
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('D30').Value = '4.244'
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').Value = '4.133'
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.05441'
$ws.Range('E32').Value = '  +6.90%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').Value = '1.197'
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('D34').Value = '1.860'
$ws.Range('E34').Value = '  -1.73%  '
$ws.Range('D35').Value = '0.7478'
$ws.Range('E35').Value = '  -3.32%  '
$ws.Range('E36').Value = '  -2.19%  '
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').Value = '1.303.57'
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('D39').Value = '0.01811'
$ws.Range('E39').Value = '  -3.32%  '
$ws.Range('D40').Value = '2.708'
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('D41').Value = '0.9348'
$ws.Range('E41').Value = '  -3.36%  '
$ws.Range('D42').Value = '6.093'
$ws.Range('E42').Value = '  +4.89%  '
$ws.Range('D43').Value = '104.47'
$ws.Range('E43').Value = '  -2.20%  '
$ws.Range('D44').Value = '0.9991'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '0.08031'
$ws.Range('E45').Value = '  +25.81%  '
$ws.Range('D46').Value = '1.975.93'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('D47').Value = '0.5174'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.00000000121'
$ws.Range('E48').Value = '  -2.59%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '9.447'
$ws.Range('E49').Value = '  -3.35%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '1.760'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '63.64'
$ws.Range('E51').Value = '  +0.06%  '

# Restore the normal (default) style now that the text values are locked
# in, so we don't leave a stray "Text" number format on these cells.
$dataRange.Style = "Normal"
